# A new weekly price-report row is inserted at row 52 (Fruta / hortaliza, semanal),
# pushing the existing rows 52:86 down to 53:87 and growing the used range to A1:R87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 52; everything below (old rows 52-86) shifts to 53-87.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record.
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value2 = 44893
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 100112022
$ws.Range("G52").Value = "Arveja Verde"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 120
$ws.Range("K52").Value = 20000
$ws.Range("L52").Value = 21000
$ws.Range("M52").Value = 20500
$ws.Range("N52").Value = "$/saco 25 kilos"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 820
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"
